$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CSR text in B2 (regenerated key/signature, same subject except "Cong Ty ABC")
$ws.Range("B2").Value = '-----BEGIN CERTIFICATE REQUEST-----
MIICtDCCAZwCAQAwcTEVMBMGA1UEAwwMTmd1eWVuIFZhbiBCMRQwEgYDVQQLDAtD
b25nIFR5IEFCQzEUMBIGA1UECgwLQ29uZyBUeSBBQkMxCzAJBgNVBAYTAlZOMQ8w
DQYDVQQHDAZRdWFuIDExDjAMBgNVBAgMBVRQSENNMIIBIjANBgkqhkiG9w0BAQEF
AAOCAQ8AMIIBCgKCAQEAhii3SsTFy6g/zPQpufFMbpt33NHlshemGxGdYMfwmUrH
F5WQ85Y7yhBllZfyRo6IuR+f3QVAhmACynd5QZ2pWZr7wGuFAwfuNo/egO9ctpkX
gnggOgW/n0lojVThMp/itcSRM9SVwyBvffPmG2CVawXFeV3G+cph5I509cyqxOyf
oKuwQ48da4YMYFIi+zGgeo2ilBrxL6zPK50Npmstg6MfHNvS98nOgH1PeBgo1OU7
HP8TVvDh+PfZV8AQQsY/9eGNMCBwoxiphPPj4v+dVqjwihP7UBr13+umnca0no0q
bwYBGRqtiWVcFYcndWTen/cPppflOfO2wsqgebz+GwIDAQABMA0GCSqGSIb3DQEB
BQUAA4IBAQBDjEhqfacS34Ec/C9eR7IpNljRQz6WRqPoTq03tQf4RJ9mKq+X5ku8
rd50XKvb82Lf7uiOiRHwNRzn/yTkOk8n5nl28kanNDbZf4ZVYR+bIYiBXTLumWN+
9hlmnLuO9bxHo0h1mn0WATvd0ReBE2CkGC6v2+k8Fwv44GBNUa8X+ymc/GsSyV+2
tdsNh671RsgR0RmIyTdaWWqkU4V/XZlQCOxhX6CPWuXXbMXPZ01XRsv+zls8uk3q
jkexOoUmZWN2+X7XwuWkhfti/fDG+I7fEG77sT0y7MEKULvmPuS+BaEUs21GGs+W
1s83sHwB3pzwoSww9yCZ/177PmhU7Cyy
-----END CERTIFICATE REQUEST-----
'

# Update the Certificate text in C2 (new serial number, validity dates 09:22:56 instead of
# 07:49:46, and subject CN/OU changed from "Test 8" to "Test 9" / MST ...8008 to ...9009,
# plus regenerated key/signature material)
$ws.Range("C2").Value = '"-----BEGIN CERTIFICATE-----
MIIGZDCCBEygAwIBAgIMXUQx+s7YNyZcrOJQMA0GCSqGSIb3DQEBCwUAMIHIMQsw
CQYDVQQGEwJWTjEUMBIGA1UECBMLSG8gQ2hpIE1pbmgxFDASBgNVBAcTC0hvIENo
aSBNaW5oMUAwPgYDVQQKEzdNb2JpbGUtSUQgVGVjaG5vbG9naWVzIGFuZCBTZXJ2
aWNlcyBKb2ludCBTdG9jayBDb21wYW55MScwJQYDVQQLEx5Nb2JpbGUtSUQgVGVj
aG5pY2FsIERlcGFydG1lbnQxIjAgBgNVBAMTGU1vYmlsZS1JRCBUcnVzdGVkIE5l
dHdvcmswHhcNMjIwODA4MDkyMjU2WhcNMjMwODEzMDkyMjU2WjCBkTELMAkGA1UE
BhMCVk4xEjAQBgNVBAgMCUjDoCBO4buZaTEaMBgGA1UECgwRTmd1eWVuIFZhbiBU
ZXN0IDkxFTATBgNVBAsMDENvbmcgVHkgQUJDRDEaMBgGA1UEAwwRTmd1eWVuIFZh
biBUZXN0IDkxHzAdBgoJkiaJk/IsZAEBDA9NU1Q6MTIzNDU2NzgwMDkwggEiMA0G
CSqGSIb3DQEBAQUAA4IBDwAwggEKAoIBAQCGKLdKxMXLqD/M9Cm58Uxum3fc0eWy
F6YbEZ1gx/CZSscXlZDzljvKEGWVl/JGjoi5H5/dBUCGYALKd3lBnalZmvvAa4UD
B+42j96A71y2mReCeCA6Bb+fSWiNVOEyn+K1xJEz1JXDIG998+YbYJVrBcV5Xcb5
ymHkjnT1zKrE7J+gq7BDjx1rhgxgUiL7MaB6jaKUGvEvrM8rnQ2may2Dox8c29L3
yc6AfU94GCjU5Tsc/xNW8OH499lXwBBCxj/14Y0wIHCjGKmE8+Pi/51WqPCKE/tQ
GvXf66adxrSejSpvBgEZGq2JZVwVhyd1ZN6f9w+ml+U587bCyqB5vP4bAgMBAAGj
ggGBMIIBfTAMBgNVHRMBAf8EAjAAMB8GA1UdIwQYMBaAFPNkMn2yPF3lLuBJfLTq
YhWUeC6rMHIGCCsGAQUFBwEBBGYwZDAyBggrBgEFBQcwAoYmaHR0cHM6Ly9tb2Jp
bGUtaWQudm4vcGtpL21vYmlsZS1pZC5jcnQwLgYIKwYBBQUHMAGGImh0dHA6Ly9t
b2JpbGUtaWQudm4vb2NzcC9yZXNwb25kZXIwRQYDVR0gBD4wPDA6BgsrBgEEAYHt
AwEEATArMCkGCCsGAQUFBwIBFh1odHRwczovL21vYmlsZS1pZC52bi9jcHMuaHRt
bDA0BgNVHSUELTArBggrBgEFBQcDAgYIKwYBBQUHAwQGCisGAQQBgjcKAwwGCSqG
SIb3LwEBBTAsBgNVHR8EJTAjMCGgH6AdhhtodHRwOi8vbW9iaWxlLWlkLnZuL2Ny
bC9nZXQwHQYDVR0OBBYEFPDsMM/q5ErwDwyRX3ZL5W4SNZ4UMA4GA1UdDwEB/wQE
AwIE8DANBgkqhkiG9w0BAQsFAAOCAgEAh6ujo0JZBb3ys59wsN7SC+uyxZuix7mx
Uwa0yiC0Dbp6KsTIxJ3UZirgsHkwOXiYU3i/LVZCG7X3YVcRxL6Ej9Ho2OKZzd3Q
TAdTi262Sm2LDdXTsFy/wmmkVr9v4zJgA6/KxamaroJG6Gvxk9iJs+Zwp1N91iBZ
K3JnhK7vmiy73qo4TbFxMnYe9ZNOCI6HDArd0b2GxUspjIIeUx3X9OU4VgnK8JBe
IJX9Fv+ZhhW12bns3LTNBXTrq5Oa7vKChxJtdBCpU73P5bKhr5k4/rHvAmM4Vzo1
+jJHA0g0OZOPRdUcLs8Z5F/6YTD1M41AY9Wnb+dW9XeTcJngmCP5r/1z6USl4vpz
V0td/ggIF7IRKl6q24GXFstccMj3A7egfnv0BWVbdlQJKF35rg/GtiXJ0BksyJXN
DxTmxX1RmbvK3KZRocjYlPlW88Z1wbi5MF//wbHESrGBVe2IQ9XlCNKtUhV6nSYx
IanLtwRe/TsR51jcANW8ZkgpuBMj3wOLNj/nFBM7e0RiEK8LNEKoOIL5ca05bt9b
CCYpTZ0ggJenmaSYMP58WZB+EOs97a1yn3cCiouFbGusEaG44FJVbQ6A9kzVBgfs
C9oVfwvign9orWBiIZArLAj5mAeCfa6QWRFumq9nnXkXr1mephnSarA3yjSnnOXi
3c3pDFm7LcM=
-----END CERTIFICATE-----"'

# Update the view: scroll/selection moved from I2 to J2, with row 2 pinned at the top
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J2").Select()
